$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.04
$ws.Range("H2").Value = 1.31
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1.31
$ws.Range("K2").Value = 1000
$ws.Range("P2").Value = 1.24
$ws.Range("Q2").Value = 1.01

# Row 3
$ws.Range("F3").Value = 1.04
$ws.Range("H3").Value = 1.04
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 1.01
$ws.Range("P3").Value = 1.24

# Row 4
$ws.Range("F4").Value = 3.45
$ws.Range("H4").Value = 2.34
$ws.Range("K4").Value = 3.25
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 1.32
$ws.Range("O4").Value = 1.01
$ws.Range("P4").Value = 1.32
$ws.Range("Q4").Value = 1.02
$ws.Range("S4").Value = 1.01
$ws.Range("T4").Value = 1.01
$ws.Range("U4").Value = 1.01

# Row 5
$ws.Range("F5").Value = 1.04
$ws.Range("G5").Value = 1000
$ws.Range("H5").Value = 1.04
$ws.Range("J5").Value = 1.01
$ws.Range("P5").Value = 1.24
$ws.Range("Q5").Value = 1.02

$wb.Save()
